$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column D width ---
$ws.Columns.Item(4).ColumnWidth = 42.6

# --- Row 9: blank spacer row with a bit of extra height ---
$ws.Rows.Item(9).RowHeight = 15.75

# --- Row 10: Packet duration (ms), styled like "Good" + a double border ---
$a10 = $ws.Cells.Item(10,1)
$a10.Value = "Packet duration"
$a10.Style = "Good"
$a10border = $a10.Borders
$a10border.LineStyle = -4119
$a10border.Color = 4144959

$b10 = $ws.Cells.Item(10,2)
$b10.Value = 36
$b10.Style = "Good"
$b10border = $b10.Borders
$b10border.LineStyle = -4119
$b10border.Color = 4144959

$ws.Cells.Item(10,3).Value = "ms"
$ws.Rows.Item(10).RowHeight = 16.5

# --- Row 11: Packet duration (tics), plain formatting ---
$ws.Cells.Item(11,1).Value = "Packet duration"
$ws.Cells.Item(11,2).Value = 64
$ws.Cells.Item(11,3).Value = "tics"
$ws.Rows.Item(11).RowHeight = 15.75

# --- Row 12: Address, "Good" style ---
$a12 = $ws.Cells.Item(12,1)
$a12.Value = "Address"
$a12.Style = "Good"

$b12 = $ws.Cells.Item(12,2)
$b12.Value = 10
$b12.Style = "Good"

# --- Row 13: WhenTransmit (tics), "Calculation" style, formula ---
$a13 = $ws.Cells.Item(13,1)
$a13.Value = "WhenTransmit"
$a13.Style = "Calculation"

$b13 = $ws.Cells.Item(13,2)
$b13.Formula = "=B12*B11"
$b13.Style = "Calculation"

$ws.Cells.Item(13,3).Value = "tics"

# --- Row 14: WhenTransmit (tics), plain formatting, formula ---
$ws.Cells.Item(14,1).Value = "WhenTransmit"
$ws.Cells.Item(14,2).Formula = "=B13/B6"
$ws.Cells.Item(14,3).Value = "tics"

# --- D11: set last, so "To ease multiplication" lands after Address/WhenTransmit
#     in the shared-string table (matches original authoring order) ---
$ws.Cells.Item(11,4).Value = "To ease multiplication"

# --- Selection moves to A14 ---
$ws.Range("A14").Select() | Out-Null
